$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.456.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.45%  "
$ws.Range("D3").Value = "'1.804.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.14%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "'307.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.40%  "
$ws.Range("D7").Value = "'0.4537"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.95%  "
$ws.Range("D8").Value = "'0.3648"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.84%  "
$ws.Range("D9").Value = "'0.07096"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.11%  "
$ws.Range("D10").Value = "'0.8707"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("D11").Value = "'0.07769"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").Value = "'19.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.26%  "
$ws.Range("D13").Value = "'1.832.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").Value = "'5.261"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.66%  "
$ws.Range("D15").Value = "'6.326"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.52%  "
$ws.Range("E16").Value = "  -5.96%  "
$ws.Range("D17").Value = "'1.008"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("E18").Value = "  -4.70%  "
$ws.Range("D19").Value = "'1.008"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").Value = "'26.488.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.35%  "
$ws.Range("D21").Value = "'14.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.20%  "
$ws.Range("D22").Value = "'4.947"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.67%  "
$ws.Range("D23").Value = "'2.071.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.90%  "
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("D25").Value = "'1.968"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("D26").Value = "'150.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.07%  "
$ws.Range("D27").Value = "'17.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.66%  "
$ws.Range("D28").Value = "'1.990"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.27%  "
$ws.Range("D29").Value = "'113.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.69%  "
$ws.Range("D30").Value = "'4.860"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.85%  "
$ws.Range("D31").Value = "'0.08668"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.08%  "
$ws.Range("D32").Value = "'3.107"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").Value = "'0.7265"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.89%  "
$ws.Range("D34").Value = "'4.423"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("D35").Value = "'1.111"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.46%  "
$ws.Range("D36").Value = "'1.007"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("D37").Value = "'2.502"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.59%  "
$ws.Range("D38").Value = "'1.074"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("D39").Value = "'0.01904"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.90%  "
$ws.Range("D40").Value = "'0.05074"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.31%  "
$ws.Range("D41").Value = "'2.863"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.02%  "
$ws.Range("D42").Value = "'6.876"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("D43").Value = "'0.4899"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.90%  "
$ws.Range("E44").Value = "  -4.67%  "
$ws.Range("D45").Value = "'8.119"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.58%  "
$ws.Range("D46").Value = "'1.007"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("E47").Value = "  -4.80%  "
$ws.Range("D48").Value = "'101.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("D49").Value = "'9.904"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.27%  "
$ws.Range("E50").Value = "  -4.70%  "
$ws.Range("D51").Value = "'0.05990"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.65%  "
